# Workbook edit: update the "About" sheet explanatory note and make the
# "About" sheet the active/selected sheet when the workbook is opened
# (matches eps-1.4.3-california-wipF update from Chris, 11/14/20).

$wb = $excel.ActiveWorkbook

# 1. Update the explanatory note text on the "About" sheet.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A3").Value = "This variable is not avaiable in the initial specification of the California model"

# 2. Make "About" the active sheet/tab (previously "PPRiEYFUfICaWHR" was
#    the active sheet), with cell A4 selected.
$wsAbout.Activate() | Out-Null
$wsAbout.Range("A4").Select() | Out-Null
